$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 271, shifting IAD..YHZ down by one.
$ws.Rows.Item(271).Insert()

# Copy formatting (including the bold/border style on column A) from the
# row that was just pushed down (now row 272) onto the newly blank row 271.
$ws.Range("A272:H272").Copy()
$ws.Range("A271:H271").PasteSpecial(-4122)

# Populate the new row with the Bishkek, Kyrgyzstan colo entry.
$ws.Cells.Item(271, 1).Value = "FRU"
$ws.Cells.Item(271, 2).Value = "Bishkek, Kyrgyzstan"
$ws.Cells.Item(271, 3).Value = "Asia Pacific"
$ws.Cells.Item(271, 4).Value = "Bishkek"
$ws.Cells.Item(271, 5).Value = "Kyrgyzstan"
$ws.Cells.Item(271, 6).Value = "KG"
$ws.Cells.Item(271, 7).Value = 42.875608
$ws.Cells.Item(271, 8).Value = 74.604613
